# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country name pairs in column A (order changed in the source list) ---
# Row 201 <-> Row 202 : Laos / Santa Lucia
$ws.Range("A201").Value = "Santa Lucia"
$ws.Range("A202").Value = "Laos"

# Row 203 <-> Row 204 : Dominica / Fiyi
$ws.Range("A203").Value = "Fiyi"
$ws.Range("A204").Value = "Dominica"

# Row 208 <-> Row 209 : Islas Malvinas / Groenlandia
$ws.Range("A208").Value = "Groenlandia"
$ws.Range("A209").Value = "Islas Malvinas"

# Row 212 <-> Row 213 : Montserrat / Seychelles
$ws.Range("A212").Value = "Seychelles"
$ws.Range("A213").Value = "Montserrat"

# --- Update the "Datos actualizados" timestamp text (cell A1) ---
$ws.Range("A1").Value = "Datos actualizados a 27 de Junio de 2020 a las 07:15"

# --- Update statistics for Pakistan (row 15) ---
$ws.Range("B15").Value = 198883
$ws.Range("C15").Value = 3138
$ws.Range("D15").Value = 86906
$ws.Range("E15").Value = 107942
$ws.Range("G15").Value = 73
$ws.Range("H15").Value = 4035

# --- Update statistics for Australia (row 74) ---
$ws.Range("B74").Value = 7641
$ws.Range("C74").Value = 46
$ws.Range("D74").Value = 6979
$ws.Range("E74").Value = 558

# --- Update statistics for Uzbekistan (row 75) ---
$ws.Range("B75").Value = 7490
$ws.Range("C75").Value = 63
$ws.Range("E75").Value = 2432

# --- Update statistics for Kirguistan (row 90) ---
$ws.Range("B90").Value = 4446
$ws.Range("C90").Value = 242
$ws.Range("D90").Value = 2194
$ws.Range("E90").Value = 2206
$ws.Range("G90").Value = 3
$ws.Range("H90").Value = 46

# --- Update statistics for Tailandia (row 98) ---
$ws.Range("D98").Value = 3053
$ws.Range("E98").Value = 51

# --- Update statistics for row 212 (now Seychelles) ---
$ws.Range("D212").Value = 11
$ws.Range("H212").Value = 0

# --- Update statistics for row 213 (now Montserrat) ---
$ws.Range("H213").Value = 1
